$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 1370.4166
$ws.Range("I2").Value = 374.33334
$ws.Range("J2").Value = 2366.5
$ws.Range("K2").Value = 374.33334
$ws.Range("L2").Value = 2366.5
$ws.Range("M2").Value = -261.33334
$ws.Range("N2").Value = -2592.5

$ws.Range("H11").Value = 0
$ws.Range("I11").Value = 0
$ws.Range("K11").Value = 0
$ws.Range("M11").ClearContents()

$ws.Range("H28").Value = 1148.1177
$ws.Range("I28").Value = 594.6
$ws.Range("K28").Value = 594.6
$ws.Range("M28").Value = -109.6

$ws.Range("H42").Value = 511.125
$ws.Range("J42").Value = 231.66667
$ws.Range("L42").Value = 695.00001
$ws.Range("N42").Value = -1155.00001

$ws.Range("H51").Value = 15629920
$ws.Range("J51").Value = 8848.5
$ws.Range("L51").Value = 8848.5
$ws.Range("N51").Value = -9816.5

$ws.Range("H132").Value = 2564740.8
$ws.Range("I132").Value = 2659677
$ws.Range("J132").Value = 1464
$ws.Range("K132").Value = 7979031
$ws.Range("L132").Value = 4392
$ws.Range("M132").Value = -7976501
$ws.Range("N132").Value = -9452

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15696.358
$ws.Range("I32").Value = 15696.358
$ws.Range("K32").Value = 15696.358
$ws.Range("M32").Value = -15409.358

$ws.Range("H110").Value = 1395.9
$ws.Range("I110").Value = 1370.5
$ws.Range("K110").Value = 1370.5
$ws.Range("M110").Value = 674.5

$ws.Range("H124").Value = 54399.8
$ws.Range("J124").Value = 54399.8
$ws.Range("L124").Value = 54399.8
$ws.Range("N124").Value = -64219.8

$ws.Range("H132").Value = 1721.7142
$ws.Range("I132").Value = 1258.4524
$ws.Range("J132").Value = 2648.238
$ws.Range("K132").Value = 3775.357199999999
$ws.Range("L132").Value = 7944.714
$ws.Range("M132").Value = -1245.357199999999
$ws.Range("N132").Value = -13004.714

$ws.Range("H137").Value = 107288.78
$ws.Range("J137").Value = 107288.78
$ws.Range("L137").Value = 107288.78
$ws.Range("N137").Value = -117488.78

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 2255.4849
$ws.Range("I99").Value = 1037.7273
$ws.Range("K99").Value = 1037.7273
$ws.Range("M99").Value = 460.2727

$ws.Range("H107").Value = 35417.133
$ws.Range("J107").Value = 2056.5715
$ws.Range("L107").Value = 2056.5715
$ws.Range("N107").Value = -5896.5715

$ws.Range("H133").Value = 99780
$ws.Range("J133").Value = 99780
$ws.Range("L133").Value = 99780
$ws.Range("N133").Value = -109900

$ws.Range("H134").Value = 2362.282
$ws.Range("I134").Value = 2447
$ws.Range("K134").Value = 7341
$ws.Range("M134").Value = -4806

$ws.Range("H135").Value = 89408.89
$ws.Range("J135").Value = 89408.89
$ws.Range("L135").Value = 89408.89
$ws.Range("N135").Value = -99548.89

$ws.Range("H137").Value = 193163
$ws.Range("I137").Value = 180709
$ws.Range("J137").Value = 199390
$ws.Range("K137").Value = 180709
$ws.Range("L137").Value = 199390
$ws.Range("M137").Value = -175609
$ws.Range("N137").Value = -209590

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H112").Value = 0
$ws.Range("J112").Value = 0
$ws.Range("L112").Value = 0
$ws.Range("N112").ClearContents()

$ws.Range("H135").Value = 111952.18
$ws.Range("J135").Value = 111952.18
$ws.Range("L135").Value = 111952.18
$ws.Range("N135").Value = -122092.18

$ws.Range("H137").Value = 84267.664
$ws.Range("J137").Value = 84267.664
$ws.Range("L137").Value = 84267.664
$ws.Range("N137").Value = -94467.664

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 47652224
$ws.Range("I4").Value = 31577636
$ws.Range("K4").Value = 94732908
$ws.Range("M4").Value = -94732796

$ws.Range("H40").Value = 259.57895
$ws.Range("I40").Value = 107.375
$ws.Range("K40").Value = 429.5
$ws.Range("M40").Value = -360.5

$ws.Range("H97").Value = 1908.8572
$ws.Range("I97").Value = 1074.4
$ws.Range("J97").Value = 3995
$ws.Range("K97").Value = 3223.2
$ws.Range("L97").Value = 11985
$ws.Range("M97").Value = -2727.2
$ws.Range("N97").Value = -12977

$ws.Range("H107").Value = 1322.05
$ws.Range("I107").Value = 1558.8
$ws.Range("J107").Value = 1085.3
$ws.Range("K107").Value = 4676.4
$ws.Range("L107").Value = 3255.9
$ws.Range("M107").Value = -2756.4
$ws.Range("N107").Value = -7095.9

$ws.Range("H134").Value = 1459.4
$ws.Range("I134").Value = 1459.4
$ws.Range("K134").Value = 4378.200000000001
$ws.Range("M134").Value = 691.7999999999993

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 3779049
$ws.Range("J44").Value = 4314151.5
$ws.Range("L44").Value = 4314151.5
$ws.Range("N44").Value = -4315343.5

$ws.Range("H64").Value = 20000
$ws.Range("I64").Value = 20000
$ws.Range("K64").Value = 20000
$ws.Range("M64").Value = -19752

$ws.Range("H67").Value = 20000
$ws.Range("I67").Value = 20000
$ws.Range("K67").Value = 20000
$ws.Range("M67").Value = -19142

$ws.Range("H70").Value = 5761.2
$ws.Range("I70").Value = 5393.7827
$ws.Range("K70").Value = 5393.7827
$ws.Range("M70").Value = -5123.7827

$ws.Range("H73").Value = 5761.2
$ws.Range("I73").Value = 5393.7827
$ws.Range("K73").Value = 5393.7827
$ws.Range("M73").Value = -4457.7827

$ws.Range("H113").Value = 1892.7
$ws.Range("I113").Value = 1899.625
$ws.Range("J113").Value = 1865
$ws.Range("K113").Value = 1899.625
$ws.Range("L113").Value = 1865
$ws.Range("M113").Value = 270.375
$ws.Range("N113").Value = -6205

$ws.Range("H132").Value = 675.5238000000001
$ws.Range("J132").Value = 536.6
$ws.Range("L132").Value = 1609.8
$ws.Range("N132").Value = -6669.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 5483.625
$ws.Range("I16").Value = 5997.5
$ws.Range("J16").Value = 4969.75
$ws.Range("K16").Value = 5997.5
$ws.Range("L16").Value = 4969.75
$ws.Range("M16").Value = -5827.5
$ws.Range("N16").Value = -5309.75

$ws.Range("H40").Value = 2181
$ws.Range("I40").Value = 2172.8948
$ws.Range("K40").Value = 2172.8948
$ws.Range("M40").Value = -2036.8948

$ws.Range("H46").Value = 3723.5386
$ws.Range("I46").Value = 810.4167
$ws.Range("J46").Value = 6220.5
$ws.Range("K46").Value = 810.4167
$ws.Range("L46").Value = 6220.5
$ws.Range("M46").Value = -622.4167
$ws.Range("N46").Value = -6596.5

$ws.Range("H100").Value = 2632.111
$ws.Range("I100").Value = 2399.8
$ws.Range("K100").Value = 2399.8
$ws.Range("M100").Value = -1858.8

$ws.Range("H122").Value = 10044.654
$ws.Range("I122").Value = 12422.643
$ws.Range("J122").Value = 7270.3335
$ws.Range("K122").Value = 37267.929
$ws.Range("L122").Value = 21811.0005
$ws.Range("M122").Value = -34817.929
$ws.Range("N122").Value = -26711.0005

$ws.Range("H132").Value = 3546.697
$ws.Range("I132").Value = 3011.25
$ws.Range("K132").Value = 9033.75
$ws.Range("M132").Value = -6503.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 5742.9165
$ws.Range("I62").Value = 3601.8572
$ws.Range("K62").Value = 3601.8572
$ws.Range("M62").Value = -2977.8572

$ws.Range("H65").Value = 5742.9165
$ws.Range("I65").Value = 3601.8572
$ws.Range("K65").Value = 18009.286
$ws.Range("M65").Value = -14889.286

$ws.Range("H81").Value = 5063.615
$ws.Range("J81").Value = 3616.75
$ws.Range("L81").Value = 7233.5
$ws.Range("N81").Value = -9355.5

$ws.Range("H84").Value = 5063.615
$ws.Range("J84").Value = 3616.75
$ws.Range("L84").Value = 36167.5
$ws.Range("N84").Value = -46775.5

$ws.Range("H122").Value = 84355.75
$ws.Range("I122").Value = 97281.71000000001
$ws.Range("K122").Value = 291845.13
$ws.Range("M122").Value = -289395.13

$ws.Range("H126").Value = 265550
$ws.Range("I126").Value = 2163.4
$ws.Range("K126").Value = 6490.200000000001
$ws.Range("M126").Value = -4020.200000000001

$ws.Range("H132").Value = 24263.824
$ws.Range("I132").Value = 27894.035
$ws.Range("J132").Value = 3208.6
$ws.Range("K132").Value = 83682.105
$ws.Range("L132").Value = 9625.799999999999
$ws.Range("M132").Value = -81152.105
$ws.Range("N132").Value = -14685.8
